$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-18 Thursday", "2025-12-19 Friday"),
    @("587÷8=73, 3", "819÷7=117, 0"),
    @("835÷6=139, 1", "445÷9=49, 4"),
    @("498÷8=62, 2", "576÷4=144, 0"),
    @("150÷3=50, 0", "785÷2=392, 1"),
    @("589÷5=117, 4", "975÷5=195, 0"),
    @("400÷9=44, 4", "434÷3=144, 2"),
    @("123÷8=15, 3", "855÷7=122, 1"),
    @("235÷5=47, 0", "475÷9=52, 7"),
    @("310÷2=155, 0", "872÷7=124, 4"),
    @("583÷7=83, 2", "946÷9=105, 1"),
    @("516÷5=103, 1", "498÷6=83, 0"),
    @("151÷6=25, 1", "500÷9=55, 5"),
    @("482÷3=160, 2", "185÷3=61, 2"),
    @("102÷5=20, 2", "186÷2=93, 0"),
    @("420÷2=210, 0", "672÷9=74, 6"),
    @("229÷6=38, 1", "437÷7=62, 3"),
    @("415÷4=103, 3", "290÷5=58, 0"),
    @("191÷9=21, 2", "471÷7=67, 2"),
    @("207÷3=69, 0", "633÷6=105, 3"),
    @("726÷5=145, 1", "154÷6=25, 4"),
    @("839÷9=93, 2", "750÷2=375, 0"),
    @("345÷7=49, 2", "824÷9=91, 5"),
    @("562÷9=62, 4", "844÷2=422, 0"),
    @("779÷2=389, 1", "360÷7=51, 3"),
    @("631÷9=70, 1", "321÷3=107, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
